$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2347
$ws.Range("E2").Value = 67
$ws.Range("F2").Value = 67
$ws.Range("G2").Value = 73
$ws.Range("H2").Value = 35
$ws.Range("I2").Value = 35
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2378
$ws.Range("L2").Value = 785
$ws.Range("M2").Value = 1593
$ws.Range("N2").Value = 1593
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 59
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = -142
$ws.Range("S2").Value = 124
$ws.Range("T2").Value = 120
$ws.Range("U2").Value = -111
$ws.Range("V2").Value = 242
$ws.Range("W2").Value = 2.87
$ws.Range("X2").Value = 1.48
$ws.Range("Y2").Value = 2.2
$ws.Range("Z2").Value = 1.47
$ws.Range("AA2").Value = 49.23
$ws.Range("AB2").Value = 2701.31
$ws.Range("AC2").Value = 294
$ws.Range("AD2").Value = 27.99
$ws.Range("AE2").Value = 13469
$ws.Range("AF2").Value = 0.61
$ws.Range("AG2").Value = 80
$ws.Range("AH2").Value = 0.97
$ws.Range("AI2").Value = 27.18
$ws.Range("AJ2").Value = 11828858

# Row 3
$ws.Range("D3").Value = 2286
$ws.Range("E3").Value = 36
$ws.Range("F3").Value = 36
$ws.Range("G3").Value = 31
$ws.Range("H3").Value = 21
$ws.Range("I3").Value = 21
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2272
$ws.Range("L3").Value = 672
$ws.Range("M3").Value = 1600
$ws.Range("N3").Value = 1600
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 59
$ws.Range("Q3").Value = 161
$ws.Range("R3").Value = -59
$ws.Range("S3").Value = -34
$ws.Range("T3").Value = 77
$ws.Range("U3").Value = 84
$ws.Range("V3").Value = 221
$ws.Range("W3").Value = 1.57
$ws.Range("X3").Value = 0.9
$ws.Range("Y3").Value = 1.29
$ws.Range("Z3").Value = 0.88
$ws.Range("AA3").Value = 41.98
$ws.Range("AB3").Value = 2722.88
$ws.Range("AC3").Value = 174
$ws.Range("AD3").Value = 52.63
$ws.Range("AE3").Value = 13525
$ws.Range("AF3").Value = 0.68
$ws.Range("AG3").Value = 80
$ws.Range("AH3").Value = 0.87
$ws.Range("AI3").Value = 46.01
$ws.Range("AJ3").Value = 11828858

# Row 4
$ws.Range("D4").Value = 2069
$ws.Range("E4").Value = -81
$ws.Range("F4").Value = -81
$ws.Range("G4").Value = -120
$ws.Range("H4").Value = -124
$ws.Range("I4").Value = -124
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2467
$ws.Range("L4").Value = 1002
$ws.Range("M4").Value = 1465
$ws.Range("N4").Value = 1465
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 59
$ws.Range("Q4").Value = -52
$ws.Range("R4").Value = -227
$ws.Range("S4").Value = 297
$ws.Range("T4").Value = 152
$ws.Range("U4").Value = -204
$ws.Range("V4").Value = 533
$ws.Range("W4").Value = -3.89
$ws.Range("X4").Value = -5.99
$ws.Range("Y4").Value = -8.09
$ws.Range("Z4").Value = -5.23
$ws.Range("AA4").Value = 68.43000000000001
$ws.Range("AB4").Value = 2499.76
$ws.Range("AC4").Value = -1048
$ws.Range("AD4").Value = -6.22
$ws.Range("AE4").Value = 12382
$ws.Range("AF4").Value = 0.53
$ws.Range("AG4").Value = 80
$ws.Range("AH4").Value = 1.23
$ws.Range("AI4").Value = -7.64
$ws.Range("AJ4").Value = 11828858

# Row 5
$ws.Range("D5").Value = 2865
$ws.Range("E5").Value = 38
$ws.Range("F5").Value = 38
$ws.Range("G5").Value = -33
$ws.Range("H5").Value = -82
$ws.Range("I5").Value = -82
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2598
$ws.Range("L5").Value = 1206
$ws.Range("M5").Value = 1392
$ws.Range("N5").Value = 1391
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 59
$ws.Range("Q5").Value = -158
$ws.Range("R5").Value = 18
$ws.Range("S5").Value = -2
$ws.Range("T5").Value = 110
$ws.Range("U5").Value = -268
$ws.Range("V5").Value = 527
$ws.Range("W5").Value = 1.31
$ws.Range("X5").Value = -2.88
$ws.Range("Y5").Value = -5.77
$ws.Range("Z5").Value = -3.25
$ws.Range("AA5").Value = 86.67
$ws.Range("AB5").Value = 2352.66
$ws.Range("AC5").Value = -697
$ws.Range("AD5").Value = -8.25
$ws.Range("AE5").Value = 11762
$ws.Range("AF5").Value = 0.49
$ws.Range("AG5:AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 11828858

# Row 6
$ws.Range("D6").Value = 1633
$ws.Range("E6").Value = -30
$ws.Range("F6").Value = -30
$ws.Range("G6").Value = -25
$ws.Range("H6").Value = -139
$ws.Range("I6").Value = -139
$ws.Range("K6").Value = 2053
$ws.Range("L6").Value = 802
$ws.Range("M6").Value = 1251
$ws.Range("N6").Value = 1244
$ws.Range("P6").Value = 59
$ws.Range("Q6").Value = 143
$ws.Range("R6").Value = -22
$ws.Range("S6").Value = -159
$ws.Range("T6").Value = 25
$ws.Range("U6").Value = 118
$ws.Range("V6").Value = 377
$ws.Range("W6").Value = -1.82
$ws.Range("X6").Value = -8.52
$ws.Range("Y6").Value = -10.55
$ws.Range("Z6").Value = -5.98
$ws.Range("AA6").Value = 64.14
$ws.Range("AB6").Value = 2107.32
$ws.Range("AC6").Value = -1176
$ws.Range("AD6").Value = -3.79
$ws.Range("AE6").Value = 10520
$ws.Range("AF6").Value = 0.42
$ws.Range("AG6:AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 11828858

# Rows 7-9: clear all data columns, leaving only A/B/C
$ws.Range("D7:AJ9").ClearContents()
